$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 444, shifting existing rows 444:569 down to 445:570
$ws.Rows.Item(444).Insert()

# Populate the newly inserted row 444 with the new data record
$ws.Cells.Item(444, 1).Value = 3
$ws.Cells.Item(444, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(444, 3).Value = "Coquimbo"
$ws.Cells.Item(444, 4).Value = 44841
$ws.Cells.Item(444, 5).Value = 5
$ws.Cells.Item(444, 6).Value = 100112037
$ws.Cells.Item(444, 7).Value = "Cebollín"
$ws.Cells.Item(444, 8).Value = "Sin especificar"
$ws.Cells.Item(444, 9).Value = "Primera"
$ws.Cells.Item(444, 10).Value = 185
$ws.Cells.Item(444, 11).Value = 3500
$ws.Cells.Item(444, 12).Value = 3800
$ws.Cells.Item(444, 13).Value = 3605
$ws.Cells.Item(444, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(444, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(444, 16).Value = 100
$ws.Cells.Item(444, 17).Value = 36
$ws.Cells.Item(444, 18).Value = "Hortaliza"
